# Append 15 new NBA game rows (933-947) to Sheet1, mirroring the
# box-score layout of the existing data (A:Away team, B:Away Pts,
# C:Home team, D:Home Pts, E:Overtime, F:Attend., G:Arena, H:Win, I:Loss).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$game1 = @('Brooklyn Nets', 112, 'Detroit Pistons', 118, 'No', 17832, 'Little Caesars Arena', 'Detroit Pistons', 'Brooklyn Nets')
$game2 = @('Minnesota Timberwolves', 113, 'Indiana Pacers', 111, 'No', 17832, 'Gainbridge Fieldhouse', 'Minnesota Timberwolves', 'Indiana Pacers')
$game3 = @('Miami Heat', 108, 'Dallas Mavericks', 114, 'No', 17832, 'American Airlines Center', 'Dallas Mavericks', 'Miami Heat')
$game4 = @('Toronto Raptors', 113, 'Phoenix Suns', 120, 'No', 17832, 'Footprint Center', 'Phoenix Suns', 'Toronto Raptors')
$game5 = @('Boston Celtics', 109, 'Denver Nuggets', 115, 'No', 17832, 'Ball Arena', 'Denver Nuggets', 'Boston Celtics')
$game6 = @('Chicago Bulls', 125, 'Golden State Warriors', 122, 'No', 17832, 'Chase Center', 'Chicago Bulls', 'Golden State Warriors')
$game7 = @('San Antonio Spurs', 129, 'Sacramento Kings', 131, 'No', 17832, 'Golden 1 Center', 'Sacramento Kings', 'San Antonio Spurs')
$game8 = @('New Orleans Pelicans', 103, 'Philadelphia 76ers', 95, 'No', 17832, 'Wells Fargo Center', 'New Orleans Pelicans', 'Philadelphia 76ers')
$game9 = @('Charlotte Hornets', 100, 'Washington Wizards', 112, 'No', 17832, 'Capital One Arena', 'Washington Wizards', 'Charlotte Hornets')
$game10 = @('Minnesota Timberwolves', 104, 'Cleveland Cavaliers', 113, 'OT', 17832, 'Rocket Mortgage Fieldhouse', 'Cleveland Cavaliers', 'Minnesota Timberwolves')
$game11 = @('Orlando Magic', 74, 'New York Knicks', 98, 'No', 17832, 'Madison Square Garden (IV)', 'New York Knicks', 'Orlando Magic')
$game12 = @('Atlanta Hawks', 99, 'Memphis Grizzlies', 92, 'No', 17832, 'FedEx Forum', 'Atlanta Hawks', 'Memphis Grizzlies')
$game13 = @('Miami Heat', 100, 'Oklahoma City Thunder', 107, 'No', 17832, 'Paycom Center', 'Oklahoma City Thunder', 'Miami Heat')
$game14 = @('Milwaukee Bucks', 122, 'Los Angeles Lakers', 123, 'No', 17832, 'Crypto.com Arena', 'Los Angeles Lakers', 'Milwaukee Bucks')
$game15 = @('Houston Rockets', 123, 'Portland Trail Blazers', 107, 'No', 17832, 'Moda Center', 'Houston Rockets', 'Portland Trail Blazers')

$newGames = @($game1, $game2, $game3, $game4, $game5, $game6, $game7, $game8, $game9, $game10, $game11, $game12, $game13, $game14, $game15)

$startRow = 933
for ($i = 0; $i -lt $newGames.Count; $i++) {
    $r = $startRow + $i
    $game = $newGames[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $game[$c - 1]
    }
}

$lastRow = $startRow + $newGames.Count - 1
$lastCell = $ws.Cells.Item($lastRow, 1)

# Scroll the viewport toward the newly appended rows and move the
# active selection to the final populated cell (A947).
$win = $excel.ActiveWindow
$win.ScrollRow = 916
$win.ScrollColumn = 1
$lastCell.Select()
